# Update Data Sources in Table
# - Rename existing "Sheet1" -> "Themes Page"
# - Add a new "Landing Page" sheet (placed first) summarising data source status

$wb = $excel.ActiveWorkbook

$themes = $wb.Worksheets.Item(1)
$themes.Name = "Themes Page"

$landing = $wb.Worksheets.Add($themes)
$landing.Name = "Landing Page"

# ---- Header row (row 17) ----
$landing.Range("A17").Value = "Card"
$landing.Range("B17").Value = "Full chart"
$landing.Range("C17").Value = "Source"
$landing.Range("D17").Value = "Source file"
$landing.Range("E17").Value = "Rename File"
$landing.Range("F17").Value = "Verify Source Link"
$landing.Range("G17").Value = "Current most recent"
$landing.Range("H17").Value = "Available"
$landing.Range("I17").Value = "Notes"
$landing.Range("A17:I17").Font.Bold = $true

# ---- Row 18: Annual Population ----
$landing.Range("A18").Value = "Annual Population"
$landing.Range("B18").Value = "/themes#population `n'Population of Dublin'"
$landing.Range("C18").Value = "https://www.cso.ie/px/pxeirestat/Statire/SelectVarVal/Define.asp?maintable=CNA13&PLanguage=0"
$landing.Range("D18").Value = "/data/Demographics/CNA13.csv"
$landing.Range("B18:H18").WrapText = $true

# ---- Row 19: Unemployment Quarterly Count ----
$landing.Range("A19").Value = "Unemployment Quarterly Count"
$landing.Range("B19").Value = "/themes#employment `n'Numbers Unemployed'"
$landing.Range("C19").Value = "https://www.cso.ie/px/pxeirestat/Statire/SelectVarVal/Define.asp?maintable=QLF08&PLanguage=0"
$landing.Range("D19").Value = "/data/Economy/QNQ22_2.csv"
$landing.Range("E19").Value = "/data/Economy/QL408.csv"
$landing.Range("F19").Value = "Good"
$landing.Range("G19").Value = "2018Q3"
$landing.Range("B19:G19").WrapText = $true
$landing.Range("I19").WrapText = $true
$landing.Range("H19").Value = "2019Q2"
$landing.Range("H19").WrapText = $true
$landing.Range("H19").Font.Color = 255

# ---- Row 20: Monthly Property Price Index ----
$landing.Range("A20").Value = "Monthly Property Price Index"
$landing.Range("B20").Value = "/themes#property-price-monthly `n'Monthly Monthly Residential Property Price Index'"
$landing.Range("C20").Value = "https://www.cso.ie/px/pxeirestat/Statire/SelectVarVal/Define.asp?maintable=HPM06&PLanguage=0"
$landing.Range("D20").Value = "/data/Housing/HPM06.csv"
$landing.Range("B20:H20").WrapText = $true
$landing.Range("I20").Value = "Source link is wrong"

# ---- Row 21: Monthly House Unit Completions ----
$landing.Range("A21").Value = "Monthly House Unit Completions"
$landing.Range("B21").Value = "/themes#unitscompmonth`nMonthly House Unit Completions"
$landing.Range("C21").Value = "https://www.housing.gov.ie/housing/statistics/house-building-and-private-rented/construction-activity-esb-connections"
$landing.Range("B21:H21").WrapText = $true

# ---- Row 22: trailing styled-but-empty row ----
$landing.Range("B22:H22").WrapText = $true

# ---- Row heights (wrapped text rows are taller) ----
$landing.Rows(18).RowHeight = 30
$landing.Rows(19).RowHeight = 30
$landing.Rows(20).RowHeight = 45
$landing.Rows(21).RowHeight = 30

# ---- Column widths ----
$landing.Columns(1).ColumnWidth = 27.830729
$landing.Columns(2).ColumnWidth = 34.666667
$landing.Columns(3).ColumnWidth = 52.166667
$landing.Columns(4).ColumnWidth = 26.998698
$landing.Columns(5).ColumnWidth = 21.498698
$landing.Columns(6).ColumnWidth = 14.666667
$landing.Columns(7).ColumnWidth = 16.998698
$landing.Columns(9).ColumnWidth = 35.830729

# ---- View ----
$landing.Activate()
$excel.ActiveWindow.Zoom = 110
$landing.Range("C20").Select()

Write-Host "Landing Page created"
